$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.396.99'
$ws.Range('E2').Value = '  +6.60%  '
$ws.Range('D3').Value = '1.816.99'
$ws.Range('E3').Value = '  +6.57%  '
$ws.Range('E4').Value = '  +0.32%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '343.40'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.91%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.0000'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3859'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.96%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '50.53'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.35%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3544'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +7.32%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.250'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +7.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07823'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +7.04%  '
$ws.Range('E12').Value = '  +14.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.001'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.39%  '
$ws.Range('E14').Value = '  +8.05%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.280'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +6.84%  '
$ws.Range('D16').Value = '1.813.86'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001138'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +6.48%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06774'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.62%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '87.00'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +7.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '18.05'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +12.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.604'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +9.58%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.24'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.21%  '
$ws.Range('D24').Value = '27.382.66'
$ws.Range('E24').Value = '  +6.66%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.474'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.90%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.742'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +10.92%  '
$ws.Range('B27').Value = 'ImmutableX'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.530'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +20.82%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '22.15'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +16.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '154.24'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.22%  '
$ws.Range('D30').Value = '2.017.20'
$ws.Range('E30').Value = '  +6.65%  '
$ws.Range('E31').Value = '  +7.68%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.474'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +9.18%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.136'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.86%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '13.90'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +9.70%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.08849'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.22%  '
$ws.Range('E36').Value = '  +2.29%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.697'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +7.60%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.7020'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +15.29%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06604'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.02442'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +8.33%  '
$ws.Range('E41').Value = '  +7.85%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '9.107'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +7.91%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.269'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.10%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '15.03'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +7.14%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6651'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +13.90%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.000'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.039'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.30%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.202'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +10.26%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '134.00'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07345'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.79%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '81.41'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.79%  '
